# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets
# to reflect the latest generated numbers.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("展览")
$sheet1.Range("F2").Value = 8609
$sheet1.Range("F5").Value = 89
$sheet1.Range("F6").Value = 1382
$sheet1.Range("F7").Value = 138
$sheet1.Range("F10").Value = 9397
$sheet1.Range("F14").Value = 177
$sheet1.Range("F16").Value = 6370
$sheet1.Range("F17").Value = 1064
$sheet1.Range("F18").Value = 89
$sheet1.Range("F20").Value = 135

$sheet4 = $wb.Worksheets.Item("全部类型")
$sheet4.Range("F2").Value = 8609
$sheet4.Range("F5").Value = 89
$sheet4.Range("F6").Value = 1382
$sheet4.Range("F7").Value = 138
$sheet4.Range("F12").Value = 9397
$sheet4.Range("F16").Value = 177
$sheet4.Range("F18").Value = 6370
$sheet4.Range("F19").Value = 1064
$sheet4.Range("F20").Value = 89
$sheet4.Range("F22").Value = 135
